$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.702.67"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "2.274.19"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "119.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  +4.14%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.622"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0945"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("E12").Value = "  +8.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.903"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.61%  "
$ws.Range("D16").Value = "2.614.05"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "2.271.80"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "43.614.82"
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("E19").Value = "  +2.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.05%  "
$ws.Range("E22").Value = "  -2.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.23%  "
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.11%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0919"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("E35").Value = "  +3.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +14.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0384"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("E39").Value = "  +4.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.241"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.02%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +42.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.676"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +21.84%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.98%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.98%  "
